# Project Plan Process.docx - "Jira link:" paragraph update
#
# Adds a small, grey note after "Jira link:" reading
# " (invitations are sent to all the teachers involved)" and appends a new
# blank paragraph (also styled small/grey) right after that paragraph.

$d = $word.ActiveDocument

# --- Locate the "Jira link:" paragraph without relying on a hard-coded index ---
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "Jira link:"
$find.Forward = $true
$find.Wrap = 0
$found = $find.Execute()
if (-not $found) {
    throw "Could not find the 'Jira link:' paragraph"
}

# Re-seat the match in a fresh Range (mutating Find's own range has caused
# aliasing issues with InsertXML in testing) and expand it to the whole
# paragraph, including the trailing paragraph mark.
$matchStart = $find.Parent.Start
$matchEnd = $find.Parent.End
$target = $d.Range($matchStart, $matchEnd)
$target.Expand(4)  # wdParagraph

# Grab the paragraph's existing opening-tag attributes (w14:paraId, rsids...)
# straight out of the live package so the edited paragraph keeps its
# original identity instead of getting a synthetic one.
$pkg = $d.Content.WordOpenXML
$tagPattern = '<w:p\b[^>]*>(?:(?!</w:p>).)*?Jira link:'
$openTag = '<w:p>'
if ($pkg -match $tagPattern) {
    $m = $matches[0]
    $close = $m.IndexOf('>')
    $openTag = $m.Substring(0, $close + 1)
}

$shade = '<w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/><w:sz w:val="22"/><w:szCs w:val="22"/>'

$fragment = $openTag +
    '<w:pPr><w:rPr>' + $shade + '<w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Jira link:</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr>' + $shade + '<w:lang w:val="en-US"/></w:rPr><w:t>(invitations are sent to all the teachers involved)</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:rPr>' + $shade + '<w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' +
    $fragment +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($packageXml)
